$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2026-02-14 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-15 Sunday", 2) | Out-Null

# Update each table cell by position (row-major, 5 columns) to avoid any
# cross-cell substring collisions during Find/Replace.
$t = $d.Tables(1)
$values = @(
    "98+0=",
    "89-52=",
    "0+39=",
    "79+2=",
    "26-21=",
    "27+43=",
    "57-13=",
    "84-30=",
    "54+2=",
    "45+16=",
    "91+5=",
    "49+4=",
    "33+12=",
    "58-6=",
    "5+53=",
    "21-5=",
    "40+53=",
    "38-15=",
    "43+52=",
    "90-0=",
    "54-50=",
    "76+11=",
    "32+27=",
    "46-8=",
    "82+5=",
    "4+80=",
    "13+22=",
    "23+24=",
    "4+46=",
    "78-44=",
    "26-18=",
    "37-13=",
    "78+9=",
    "28-13=",
    "46-22=",
    "31+16=",
    "94-88=",
    "49+50=",
    "85-24=",
    "37+12=",
    "92-90=",
    "6+5=",
    "26+50=",
    "25-20=",
    "27+58=",
    "85-58=",
    "73-37=",
    "74-52=",
    "77-62=",
    "80-9=",
    "75-0=",
    "52+9=",
    "20+6=",
    "67-13=",
    "51+10=",
    "27+61=",
    "33+28=",
    "57-40=",
    "45-10=",
    "97-2=",
    "1+43=",
    "12+68=",
    "1-1=",
    "2+63=",
    "99-95=",
    "62-20=",
    "31-22=",
    "44+33=",
    "64-45=",
    "7+84=",
    "19+61=",
    "97-96=",
    "8+18=",
    "67+30=",
    "24+60=",
    "23+21=",
    "22+62=",
    "97-44=",
    "19+57=",
    "98-72=",
    "5+59=",
    "20+79=",
    "68-10=",
    "9+90=",
    "74-69=",
    "2+73=",
    "7+69=",
    "27+61=",
    "34+7=",
    "86+12=",
    "51-9=",
    "77-21=",
    "97-32=",
    "25-0=",
    "98-78=",
    "24+30=",
    "27+3=",
    "41+33=",
    "81-80=",
    "40-12="
)

$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cellRange = $cell.Range
        $cellRange.MoveEnd(1, -1) | Out-Null
        $cellRange.Text = $values[$idx]
        $idx = $idx + 1
    }
}

Write-Output ("Replaced " + $idx.ToString() + " cells")